# Insert a new data row above row 466 (pushing the existing rows 466..558
# down to 467..559) and populate it with a new observation, matching the
# published diff: a new "Limache" price point dated 45015 is inserted into
# the "Femacal de La Calera - Poroto verde" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 466..558 down one row, leaving a blank (but formatted) row 466.
$ws.Rows.Item(466).Insert()

# Populate the newly inserted row 466 with the new observation.
$ws.Range("A466").Value = 3
$ws.Range("B466").Value = "Femacal de La Calera"
$ws.Range("C466").Value = "Coquimbo"
$ws.Range("D466").Value = 45015
$ws.Range("E466").Value = 5
$ws.Range("F466").Value = 100112031
$ws.Range("G466").Value = "Poroto verde"
$ws.Range("H466").Value = "Magnum"
$ws.Range("I466").Value = "Primera"
$ws.Range("J466").Value = 65
$ws.Range("K466").Value = 20000
$ws.Range("L466").Value = 21000
$ws.Range("M466").Value = 20462
$ws.Range("N466").Value = "$/malla 25 kilos"
$ws.Range("O466").Value = "Limache"
$ws.Range("P466").Value = 818
$ws.Range("Q466").Value = 25
$ws.Range("R466").Value = "Hortaliza"
